$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 583, shifting rows 583:652 down to 584:653
$ws.Rows.Item(583).Insert()

# Populate the newly inserted row 583 with values.
# Columns A-I, Q, R are copied unchanged from the (shifted) row below (now 584).
$ws.Cells.Item(583, 1).Value = 3
$ws.Cells.Item(583, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(583, 3).Value = "Coquimbo"
$ws.Cells.Item(583, 4).Value = 45212
$ws.Cells.Item(583, 5).Value = 5
$ws.Cells.Item(583, 6).Value = 100112043
$ws.Cells.Item(583, 7).Value = "Pepino ensalada"
$ws.Cells.Item(583, 8).Value = "Sin especificar"
$ws.Cells.Item(583, 9).Value = "Primera"
$ws.Cells.Item(583, 10).Value = 90
$ws.Cells.Item(583, 11).Value = 14000
$ws.Cells.Item(583, 12).Value = 15000
$ws.Cells.Item(583, 13).Value = 14556
$ws.Cells.Item(583, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(583, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(583, 16).Value = 243
$ws.Cells.Item(583, 17).Value = 60
$ws.Cells.Item(583, 18).Value = "Hortaliza"
